$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear old content in the region that changes (rows 10-24); rows 1-9 and 12 are untouched by the edit
$ws.Range("A10:C24").Clear()

# Row 24 no longer exists in the final layout; drop it so the sheet ends at row 23
$ws.Rows.Item(24).Delete()

# Row 10
$ws.Cells.Item(10, 1).Value = 'Objetivos:'
$ws.Cells.Item(10, 2).Value = '11079086 - Herlandí de Souza Andrade'
$ws.Cells.Item(10, 3).Value = '11079086 - Herlandí de Souza Andrade'
$ws.Rows.Item(10).RowHeight = 60

# Row 11
$ws.Cells.Item(11, 1).Value = 'Objectives:'
$ws.Cells.Item(11, 2).Value = 'Introduce students to the overall picture of strategic marketing management, empowering them to act in the managerial marketing process under strategic and operational perspectives.'
$ws.Cells.Item(11, 3).Value = 'Introduce students to the overall picture of strategic marketing management, empowering them to act in the managerial marketing process under strategic and operational perspectives.'
$ws.Rows.Item(11).RowHeight = 60

# Row 12
$ws.Cells.Item(12, 1).Value = 'Docentes responsáveis:'
$ws.Rows.Item(12).AutoFit()

# Row 13
$ws.Cells.Item(13, 1).Value = 'Programa resumido:'
$ws.Cells.Item(13, 2).Value = 'Semestral'
$ws.Cells.Item(13, 3).Value = 'Semestral'
$ws.Rows.Item(13).RowHeight = 60

# Row 14
$ws.Cells.Item(14, 1).Value = 'Short syllabus:'
$ws.Cells.Item(14, 2).Value = 'Strategic marketing and strategic planning oriented to the market.'
$ws.Cells.Item(14, 3).Value = 'Strategic marketing and strategic planning oriented to the market.'
$ws.Rows.Item(14).RowHeight = 60

# Row 15
$ws.Cells.Item(15, 1).Value = 'Programa:'
$ws.Cells.Item(15, 2).Value = '01/01/2021'
$ws.Cells.Item(15, 3).Value = '01/01/2021'
$ws.Rows.Item(15).RowHeight = 120

# Row 16
$ws.Cells.Item(16, 1).Value = 'Syllabus:'
$ws.Cells.Item(16, 2).Value = 'STRATEGIC MARKETING1. Evolution of the Marketing Concept and Marketing System2. Marketing, Value Concept, Market Orientation3. Market Analysis and Consumer Behavior4. Marketing Information System and Market Intelligence5. Marketing Modalities6. Fundamentals of Business Strategy and Strategic Marketing7. Strategic Management and Strategic Marketing8. Analytical Tools to Evaluate Market Opportunities9. Market Segmentation and Positioning10. The Strategic Marketing Plan11. Communication: press services, SAC''S, Ombudsman'
$ws.Cells.Item(16, 3).Value = 'STRATEGIC MARKETING1. Evolution of the Marketing Concept and Marketing System2. Marketing, Value Concept, Market Orientation3. Market Analysis and Consumer Behavior4. Marketing Information System and Market Intelligence5. Marketing Modalities6. Fundamentals of Business Strategy and Strategic Marketing7. Strategic Management and Strategic Marketing8. Analytical Tools to Evaluate Market Opportunities9. Market Segmentation and Positioning10. The Strategic Marketing Plan11. Communication: press services, SAC''S, Ombudsman'
$ws.Rows.Item(16).RowHeight = 120

# Row 17
$ws.Cells.Item(17, 1).Value = 'Avaliação:'
$ws.Rows.Item(17).AutoFit()

# Row 18
$ws.Cells.Item(18, 1).Value = 'Método:'
$ws.Cells.Item(18, 2).Value = '11079086 - Herlandí de Souza Andrade'
$ws.Cells.Item(18, 3).Value = '11079086 - Herlandí de Souza Andrade'
$ws.Rows.Item(18).RowHeight = 60

# Row 19
$ws.Cells.Item(19, 1).Value = 'Critério:'
$ws.Cells.Item(19, 2).Value = 'Aulas expositivas e dialogadas; dinâmicas, projetos e trabalhos em grupo; exercícios individuais; e, seminários, debates e palestras.'
$ws.Cells.Item(19, 3).Value = 'Aulas expositivas e dialogadas; dinâmicas, projetos e trabalhos em grupo; exercícios individuais; e, seminários, debates e palestras.'
$ws.Rows.Item(19).RowHeight = 60

# Row 20
$ws.Cells.Item(20, 1).Value = 'Norma de recuperação:'
$ws.Cells.Item(20, 2).Value = 'Média Aritmética dos Projetos, Trabalhos, Exercícios e outras atividades avaliativas realizadas no decorrer da disciplina, considerando as questões relativas às Competências (Conhecimento, Habilidade e Atitude, que incluem a presença e participação dos alunos nas aulas) desenvolvidas.'
$ws.Cells.Item(20, 3).Value = 'Média Aritmética dos Projetos, Trabalhos, Exercícios e outras atividades avaliativas realizadas no decorrer da disciplina, considerando as questões relativas às Competências (Conhecimento, Habilidade e Atitude, que incluem a presença e participação dos alunos nas aulas) desenvolvidas.'
$ws.Rows.Item(20).RowHeight = 60

# Row 21
$ws.Cells.Item(21, 1).Value = 'Bibliografia:'
$ws.Cells.Item(21, 2).Value = 'NF = (MF + PR)/2, onde MF é a média final da avaliação e PR é uma prova de recuperação.'
$ws.Cells.Item(21, 3).Value = 'NF = (MF + PR)/2, onde MF é a média final da avaliação e PR é uma prova de recuperação.'
$ws.Rows.Item(21).RowHeight = 120

# Row 22
$ws.Cells.Item(22, 1).Value = 'Requisitos:'
$ws.Rows.Item(22).AutoFit()

# Row 23
$ws.Cells.Item(23, 2).Value = 'LOQ4240 -  Administração e Organização II  (Requisito fraco)`n'
$ws.Cells.Item(23, 3).Value = 'LOQ4240 -  Administração e Organização II  (Requisito fraco)`n'
$ws.Rows.Item(23).RowHeight = 30
